$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to re-pulled/recalculated data
$ws.Range("F4").Value = 4
$ws.Range("F9").Value = -10
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 3
$ws.Range("F14").Value = 5
